$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows to insert before the existing row 231 ("Remessa", "Material", "Quantidade")
$newRows = @(
    @("80266137", "30152-OSR-I", 300000),
    @("80266148", "10000-LDG-I", 1),
    @("80266149", "10000-LDG-I", 1),
    @("80266150", "10000-LDG-I", 1),
    @("80266152", "10645-ARI-I", 1),
    @("80266153", "20098-CTY-I", 2),
    @("80266154", "10371-ARI-I", 1),
    @("80266156", "10000-LDG-I", 1),
    @("80266157", "10382-ARI-I", 1),
    @("80266158", "33532-ATE-I", 1),
    @("80266159", "10030-ARI-I", 1),
    @("80266160", "20104-CTY-I", 1),
    @("80266161", "10000-LDG-I", 1),
    @("80266162", "30259-OSR-I", 18000),
    @("80266163", "10247-ARI-I", 1),
    @("80266164", "10125-ARI-I", 1),
    @("80266165", "10065-ARI-I", 1),
    @("80266166", "10025-ARI-I", 1),
    @("80266167", "10636-ARI-I", 1),
    @("80266169", "10035-ARI-I", 1),
    @("80266170", "10255-ARI-I", 1),
    @("80266171", "10479-ARI-I", 1),
    @("80266173", "10362-ARI-I", 1),
    @("80266174", "10251-ARI-I", 1),
    @("80266175", "10000-LDG-I", 1),
    @("80266176", "10645-ARI-I", 1),
    @("80266177", "10020-ARI-I", 1),
    @("80266179", "22780-STM-I", 15000),
    @("80266180", "20953-CTY-I", 1),
    @("80266182", "12732-ROY-I", 15000),
)

$startRow = 231
$rowCount = $newRows.Length
$endRow = $startRow + $rowCount - 1

# Shift the existing rows (old 231:236) down by inserting fresh rows above them.
$ws.Rows($startRow.ToString() + ":" + $endRow.ToString()).Insert()

# Column A is all-digit order numbers ("80266137", ...); force them to stay
# text (matching the rest of the sheet) instead of being auto-parsed as
# numbers, then drop back to General so no stray number format lingers.
$colARange = $ws.Range("A" + $startRow + ":A" + $endRow)
$colARange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}

$colARange.NumberFormat = "general"

# The "Quantidade" column now reads right-aligned instead of centered.
$lastRow = 236 + $rowCount
$ws.Range("C2:C" + $lastRow).HorizontalAlignment = -4152

# Matches the saved view state: selection on F7, scrolled back to the top.
$ws.Range("F7").Select()
